$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Approved/Rejected + ReasonToReject values between row 26 and row 27.
# Row 26: was I26="Approved", J26=(empty)  -> becomes I26="Rejected", J26="Nil"
# Row 27: was I27="Rejected", J27="Nil"    -> becomes I27="Approved", J27=(empty)
$ws.Range("I26").Value = "Rejected"
$ws.Range("J26").Value = "Nil"
$ws.Range("I27").Value = "Approved"
$ws.Range("J27").ClearContents()

# Update the selection to match the new active selection (I26:J26).
$ws.Range("I26:J26").Select()
